# SM04_microstate_IDs_with_2D_depiction.xlsx
# Add a new "canonical SMILES" column (D) to the microstate table.
# For most microstates this is identical to the existing canonical isomeric
# SMILES (column C); for SM04_micro009 the canonical (non-isomeric) SMILES
# drops the stereo bond markers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("D2").Value = "canonical SMILES"

# New column width (approx. character width used by the source workbook)
$ws.Columns.Item(4).ColumnWidth = 36

# New column D values - canonical (non-isomeric) SMILES for each microstate
$ws.Range("D3").Value  = "c1ccc2c(c1)c(=[NH+]Cc3ccc(cc3)Cl)nc[nH]2"
$ws.Range("D4").Value  = "c1ccc2c(c1)c(ncn2)NCc3ccc(cc3)Cl"
$ws.Range("D5").Value  = "c1ccc2c(c1)c(=NCc3ccc(cc3)Cl)nc[nH]2"
$ws.Range("D6").Value  = "c1ccc2c(c1)c(ncn2)[N-]Cc3ccc(cc3)Cl"
$ws.Range("D7").Value  = "c1ccc2c(c1)c([nH+]cn2)NCc3ccc(cc3)Cl"
$ws.Range("D8").Value  = "c1ccc2c(c1)c(ncn2)[NH2+]Cc3ccc(cc3)Cl"
$ws.Range("D9").Value  = "c1ccc2c(c1)c(=NCc3ccc(cc3)Cl)[nH]cn2"
$ws.Range("D10").Value = "c1ccc2c(c1)c(nc[nH+]2)[NH2+]Cc3ccc(cc3)Cl"
$ws.Range("D11").Value = "c1ccc2c(c1)c([nH+]c[nH+]2)NCc3ccc(cc3)Cl"
$ws.Range("D12").Value = "c1ccc2c(c1)c([nH+]cn2)[NH2+]Cc3ccc(cc3)Cl"
$ws.Range("D13").Value = "c1ccc2c(c1)c([nH+]c[nH+]2)[NH2+]Cc3ccc(cc3)Cl"

Write-Host "Applied SAMPL6 canonical SMILES column update"
